$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.021590333333334
$ws.Range("H2").Value = 6.064771
$ws.Range("I2").Value = 0.01116262347650641
$ws.Range("J2").Value = 0.01116262347650641
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 56.98117766666667
$ws.Range("N2").Value = 170.943533
$ws.Range("O2").Value = 0.952030123851636
$ws.Range("P2").Value = 0.9520301238516359
$ws.Range("Q2").Value = 115.1925979528826
$ws.Range("R2").Value = 1036.733381575943
$ws.Range("S2").Value = 0.01062715381084758
$ws.Range("T2").Value = 0.01062715381084758

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.021590333333334
$ws.Range("H3").Value = 6.064771
$ws.Range("I3").Value = 0.01116262347650641
$ws.Range("J3").Value = 0.01116262347650641
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.516719
$ws.Range("N3").Value = 7.550157
$ws.Range("O3").Value = 0.04204883786863874
$ws.Range("P3").Value = 0.04204883786863874
$ws.Range("Q3").Value = 5.087774802116335
$ws.Range("R3").Value = 45.78997321904701
$ws.Range("S3").Value = 0.0004693753447522787
$ws.Range("T3").Value = 0.0004693753447522787

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.021590333333334
$ws.Range("H4").Value = 6.064771
$ws.Range("I4").Value = 0.01116262347650641
$ws.Range("J4").Value = 0.01116262347650641
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.3543876666666666
$ws.Range("N4").Value = 1.063163
$ws.Range("O4").Value = 0.005921038279725251
$ws.Range("P4").Value = 0.005921038279725251
$ws.Range("Q4").Value = 0.7164266811858888
$ws.Range("R4").Value = 6.447840130673
$ws.Range("S4").Value = 0.00006609432090655423
$ws.Range("T4").Value = 0.00006609432090655423

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 25.140634
$ws.Range("H5").Value = 75.421902
$ws.Range("I5").Value = 0.1388191398995883
$ws.Range("J5").Value = 0.1388191398995883
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 56.98117766666667
$ws.Range("N5").Value = 170.943533
$ws.Range("O5").Value = 0.952030123851636
$ws.Range("P5").Value = 0.9520301238516359
$ws.Range("Q5").Value = 1432.542932606641
$ws.Range("R5").Value = 12892.88639345977
$ws.Range("S5").Value = 0.1321600029515826
$ws.Range("T5").Value = 0.1321600029515826

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 25.140634
$ws.Range("H6").Value = 75.421902
$ws.Range("I6").Value = 0.1388191398995883
$ws.Range("J6").Value = 0.1388191398995883
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.516719
$ws.Range("N6").Value = 7.550157
$ws.Range("O6").Value = 0.04204883786863874
$ws.Range("P6").Value = 0.04204883786863874
$ws.Range("Q6").Value = 63.27191125984601
$ws.Range("R6").Value = 569.447201338614
$ws.Range("S6").Value = 0.005837183506701668
$ws.Range("T6").Value = 0.005837183506701667

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 25.140634
$ws.Range("H7").Value = 75.421902
$ws.Range("I7").Value = 0.1388191398995883
$ws.Range("J7").Value = 0.1388191398995883
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.3543876666666666
$ws.Range("N7").Value = 1.063163
$ws.Range("O7").Value = 0.005921038279725251
$ws.Range("P7").Value = 0.005921038279725251
$ws.Range("Q7").Value = 8.909530621780666
$ws.Range("R7").Value = 80.18577559602599
$ws.Range("S7").Value = 0.0008219534413039972
$ws.Range("T7").Value = 0.0008219534413039971

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 153.9412893333333
$ws.Range("H8").Value = 461.823868
$ws.Range("I8").Value = 0.8500182366239053
$ws.Range("J8").Value = 0.8500182366239052
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 56.98117766666667
$ws.Range("N8").Value = 170.943533
$ws.Range("O8").Value = 0.952030123851636
$ws.Range("P8").Value = 0.9520301238516359
$ws.Range("Q8").Value = 8771.755957738405
$ws.Range("R8").Value = 78945.80361964565
$ws.Range("S8").Value = 0.8092429670892057
$ws.Range("T8").Value = 0.8092429670892056

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 153.9412893333333
$ws.Range("H9").Value = 461.823868
$ws.Range("I9").Value = 0.8500182366239053
$ws.Range("J9").Value = 0.8500182366239052
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.516719
$ws.Range("N9").Value = 7.550157
$ws.Range("O9").Value = 0.04204883786863874
$ws.Range("P9").Value = 0.04204883786863874
$ws.Range("Q9").Value = 387.4269677496974
$ws.Range("R9").Value = 3486.842709747276
$ws.Range("S9").Value = 0.0357422790171848
$ws.Range("T9").Value = 0.03574227901718479

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 153.9412893333333
$ws.Range("H10").Value = 461.823868
$ws.Range("I10").Value = 0.8500182366239053
$ws.Range("J10").Value = 0.8500182366239052
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.3543876666666666
$ws.Range("N10").Value = 1.063163
$ws.Range("O10").Value = 0.005921038279725251
$ws.Range("P10").Value = 0.005921038279725251
$ws.Range("Q10").Value = 54.55489433049821
$ws.Range("R10").Value = 490.9940489744839
$ws.Range("S10").Value = 0.0050329905175147
$ws.Range("T10").Value = 0.005032990517514699
